# Daily "Updated symbol list" refresh for the crypto price sheet.
# The upstream scraper re-pulls coinranking.com data each run, so most
# rows just get a tiny price tweak in column D, a couple of rows see a
# full re-rank (coin name / link / price / rank-label shifting up one
# row, with "One" dropping out of the Top-24-worst slot and landing in
# the Top-24-best slot instead), and two "Worst/Best in 24h" labels move
# to different rows.
#
# All of the D-column prices are stored as literal text (e.g. "0.8700"
# keeps its trailing zero, "244.71" must not become a float), so each
# numeric-looking price is written with a leading apostrophe to force
# text entry, then the cell style is reset to "Normal" so we don't leave
# a stray quote-prefixed/text-formatted style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice($a1, $text) {
    $ws.Range($a1).Value = "'" + $text
    $ws.Range($a1).Style = "Normal"
}

function Set-Text($a1, $text) {
    $ws.Range($a1).Value = $text
}

# --- rows 2-16: small price updates only -----------------------------
Set-TextPrice "D2"  "244.71"
Set-TextPrice "D3"  "23.94"
Set-TextPrice "D4"  "5.198"
Set-TextPrice "D5"  "0.05732"
Set-TextPrice "D6"  "6.486"
Set-TextPrice "D7"  "3.175"
Set-TextPrice "D8"  "0.8144"
Set-TextPrice "D9"  "0.8714"
Set-TextPrice "D11" "0.06944"
Set-TextPrice "D12" "0.03156"
Set-TextPrice "D13" "0.02935"
Set-TextPrice "D14" "0.09329"
Set-TextPrice "D15" "3.854"
Set-TextPrice "D16" "0.001519"

# --- rows 18-24: re-ranked coins (name/link/price/rank-label shift) --
Set-Text      "B18" "TigerCash"
Set-Text      "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextPrice "D18" "0.006156"
Set-Text      "E18" "17TigerCashTCH"

Set-Text      "B19" "BitKan"
Set-Text      "C19" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextPrice "D19" "0.001243"
Set-Text      "E19" "18BitKanKAN"

Set-Text      "B20" "HotbitToken"
Set-Text      "C20" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextPrice "D20" "0.004106"
Set-Text      "E20" "19HotbitTokenHTB"

Set-Text      "B21" "NitroEx"
Set-Text      "C21" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextPrice "D21" "0.00008505"
Set-Text      "E21" "20NitroExNTX"

Set-Text      "B22" "LEO"
Set-Text      "C22" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextPrice "D22" "3.552"
Set-Text      "E22" "21LEOLEO"

Set-Text      "B23" "BTSEToken"
Set-Text      "C23" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextPrice "D23" "2.163"
Set-Text      "E23" "22BTSETokenBTSE"

Set-Text      "B24" "One"
Set-Text      "C24" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextPrice "D24" "0.01013"
Set-Text      "E24" "23OneONEBestin24h"

# --- remaining scattered price updates --------------------------------
Set-TextPrice "D25" "0.3192"
Set-TextPrice "D27" "0.0002331"
Set-TextPrice "D41" "0.006423"
Set-TextPrice "D42" "0.1050"

Set-TextPrice "D43" "0.002224"
Set-Text      "E43" "42CEJICEJIWorstin24h"

Set-TextPrice "D44" "0.008124"
Set-TextPrice "D45" "0.00005472"

Set-TextPrice "D48" "0.002041"
Set-Text      "E48" "47BOLOBOLO"
